$d = $word.ActiveDocument

# Fix 1: "считаетя" -> "считается" and "название песни" -> "название трека"
$d.Content.Find.Execute("БТ-5: Дубликатом считаетя файл, в котором указан тот же альбом, исполнитель и название песни, что и в другом файле.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "БТ-5: Дубликатом считается файл, в котором указан тот же альбом, исполнитель и название трека, что и в другом файле.", 2)

# Fix 2: "некорректа" -> "некорректна"
$d.Content.Find.Execute("Появляется в случае, если указанная пользователем директория некорректа.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Появляется в случае, если указанная пользователем директория некорректна.", 2)
